$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update the F column (time_taken) timestamps on the "data" sheet ---
$ws1.Cells.Item(2,6).Value = '2021-10-05 14:22:33.676556'
$ws1.Cells.Item(3,6).Value = '2021-10-05 14:22:33.676563'
$ws1.Cells.Item(4,6).Value = '2021-10-05 14:22:33.676567'
$ws1.Cells.Item(5,6).Value = '2021-10-05 14:22:33.676570'
$ws1.Cells.Item(6,6).Value = '2021-10-05 14:22:33.676573'
$ws1.Cells.Item(7,6).Value = '2021-10-05 14:22:33.676575'
$ws1.Cells.Item(8,6).Value = '2021-10-05 14:22:33.676578'
$ws1.Cells.Item(9,6).Value = '2021-10-05 14:22:33.676580'
$ws1.Cells.Item(10,6).Value = '2021-10-05 14:22:33.676583'
$ws1.Cells.Item(11,6).Value = '2021-10-05 14:22:33.676585'
$ws1.Cells.Item(12,6).Value = '2021-10-05 14:22:33.676592'
$ws1.Cells.Item(13,6).Value = '2021-10-05 14:22:33.676594'
$ws1.Cells.Item(14,6).Value = '2021-10-05 14:22:33.676597'
$ws1.Cells.Item(15,6).Value = '2021-10-05 14:22:33.676599'
$ws1.Cells.Item(16,6).Value = '2021-10-05 14:22:33.676602'
$ws1.Cells.Item(17,6).Value = '2021-10-05 14:22:33.676604'
$ws1.Cells.Item(18,6).Value = '2021-10-05 14:22:33.676607'
$ws1.Cells.Item(19,6).Value = '2021-10-05 14:22:33.676610'
$ws1.Cells.Item(20,6).Value = '2021-10-05 14:22:33.676612'
$ws1.Cells.Item(21,6).Value = '2021-10-05 14:22:33.676615'
$ws1.Cells.Item(22,6).Value = '2021-10-05 14:22:33.676617'
$ws1.Cells.Item(23,6).Value = '2021-10-05 14:22:33.676620'
$ws1.Cells.Item(24,6).Value = '2021-10-05 14:22:33.676623'
$ws1.Cells.Item(25,6).Value = '2021-10-05 14:22:33.676625'
$ws1.Cells.Item(26,6).Value = '2021-10-05 14:22:33.676628'
$ws1.Cells.Item(27,6).Value = '2021-10-05 14:22:33.676630'
$ws1.Cells.Item(28,6).Value = '2021-10-05 14:22:33.676633'
$ws1.Cells.Item(29,6).Value = '2021-10-05 14:22:33.676635'
$ws1.Cells.Item(30,6).Value = '2021-10-05 14:22:33.676638'
$ws1.Cells.Item(31,6).Value = '2021-10-05 14:22:33.676640'
$ws1.Cells.Item(32,6).Value = '2021-10-05 14:22:33.676643'
$ws1.Cells.Item(33,6).Value = '2021-10-05 14:22:33.676645'
$ws1.Cells.Item(34,6).Value = '2021-10-05 14:22:33.676648'
$ws1.Cells.Item(35,6).Value = '2021-10-05 14:22:33.676651'
$ws1.Cells.Item(36,6).Value = '2021-10-05 14:22:33.676653'
$ws1.Cells.Item(37,6).Value = '2021-10-05 14:22:33.676656'
$ws1.Cells.Item(38,6).Value = '2021-10-05 14:22:33.676658'
$ws1.Cells.Item(39,6).Value = '2021-10-05 14:22:33.676661'
$ws1.Cells.Item(40,6).Value = '2021-10-05 14:22:33.676663'
$ws1.Cells.Item(41,6).Value = '2021-10-05 14:22:33.676666'
$ws1.Cells.Item(42,6).Value = '2021-10-05 14:22:33.676668'
$ws1.Cells.Item(43,6).Value = '2021-10-05 14:22:33.676671'
$ws1.Cells.Item(44,6).Value = '2021-10-05 14:22:33.676673'
$ws1.Cells.Item(45,6).Value = '2021-10-05 14:22:33.676676'
$ws1.Cells.Item(46,6).Value = '2021-10-05 14:22:33.676679'
$ws1.Cells.Item(47,6).Value = '2021-10-05 14:22:33.676681'
$ws1.Cells.Item(48,6).Value = '2021-10-05 14:22:33.676683'
$ws1.Cells.Item(49,6).Value = '2021-10-05 14:22:33.676686'
$ws1.Cells.Item(50,6).Value = '2021-10-05 14:22:33.676688'
$ws1.Cells.Item(51,6).Value = '2021-10-05 14:22:33.676691'
$ws1.Cells.Item(52,6).Value = '2021-10-05 14:22:33.676693'
$ws1.Cells.Item(53,6).Value = '2021-10-05 14:22:33.676696'
$ws1.Cells.Item(54,6).Value = '2021-10-05 14:22:33.676698'
$ws1.Cells.Item(55,6).Value = '2021-10-05 14:22:33.676701'
$ws1.Cells.Item(56,6).Value = '2021-10-05 14:22:33.676703'
$ws1.Cells.Item(57,6).Value = '2021-10-05 14:22:33.676706'
$ws1.Cells.Item(58,6).Value = '2021-10-05 14:22:33.676709'
$ws1.Cells.Item(59,6).Value = '2021-10-05 14:22:33.676711'
$ws1.Cells.Item(60,6).Value = '2021-10-05 14:22:33.676714'
$ws1.Cells.Item(61,6).Value = '2021-10-05 14:22:33.676716'
$ws1.Cells.Item(62,6).Value = '2021-10-05 14:22:33.676718'

# --- Add the new "metadata" sheet after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Copy the bold/bordered header style from the data sheet onto the new header row
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy the index-column style onto the new sheet's A2 cell
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row ---
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# --- Data row ---
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Renal tubulopathies"
$ws2.Range("C2").Value = 292
$d2 = $ws2.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "2.28"
$ws2.Range("E2").Value = "2021-10-04T10:34:35.731460Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:33.673027"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/292/?format=json"

# Keep "data" as the active sheet, matching the original active tab
$ws1.Activate()

Write-Output "done"
